# MENU_MOCK.xlsx edit script
# 5-preset system with lower defaults, multiplier fix, xlsx generators
#
# This script updates the "Default" (column E) and "Tooltip" (column F)
# values on the Menu Options sheet to:
#  - reword preset tooltips to mention the balanced middle default
#  - rename the frequency preset default from "Normal" to "Default"
#  - fix the damage-type multiplier tooltips/defaults to describe bleed
#    multipliers (0.0x disables bleed) instead of raw damage multipliers
#  - lower several zone chance/damage/duration defaults

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: OptionDamagePreset ---
$ws.Range("F3").Value = "Damage per tick preset. Default is the balanced middle value."

# --- Row 4: OptionDurationPreset ---
$ws.Range("F4").Value = "Bleed duration preset. Default is the balanced middle value."

# --- Row 5: OptionFrequencyPreset ---
$ws.Range("E5").Value = '"Default"'
$ws.Range("F5").Value = "Tick frequency preset. Default is the balanced middle value."

# --- Row 6: OptionChancePreset ---
$ws.Range("F6").Value = "Bleed chance preset. Default is the balanced middle value."

# --- Row 7: OptionPierceMultiplier ---
$ws.Range("F7").Value = "Bleed damage multiplier for pierce attacks. 0.0x disables bleed from pierce entirely."

# --- Row 8: OptionSlashMultiplier ---
$ws.Range("F8").Value = "Bleed damage multiplier for slash attacks. 0.0x disables bleed from slash entirely."

# --- Row 9: OptionBluntMultiplier ---
$ws.Range("E9").Value = "0.5f"
$ws.Range("F9").Value = "Bleed damage multiplier for blunt attacks. 0.0x disables bleed from blunt entirely."

# --- Row 17: OptionThroatChance ---
$ws.Range("E17").Value = "60f"

# --- Row 18: OptionThroatDamage ---
$ws.Range("E18").Value = "2.5f"

# --- Row 19: OptionThroatDuration ---
$ws.Range("E19").Value = "6.0f"

# --- Row 21: OptionHeadChance ---
$ws.Range("E21").Value = "40f"

# --- Row 22: OptionHeadDamage ---
$ws.Range("E22").Value = "1.5f"

# --- Row 23: OptionHeadDuration ---
$ws.Range("E23").Value = "5.0f"

# --- Row 25: OptionNeckChance ---
$ws.Range("E25").Value = "55f"

# --- Row 26: OptionNeckDamage ---
$ws.Range("E26").Value = "2.0f"

# --- Row 27: OptionNeckDuration ---
$ws.Range("E27").Value = "5.5f"

# --- Row 29: OptionTorsoChance ---
$ws.Range("E29").Value = "35f"

# --- Row 30: OptionTorsoDamage ---
$ws.Range("E30").Value = "1.0f"

# --- Row 31: OptionTorsoDuration ---
$ws.Range("E31").Value = "4.0f"

# --- Row 33: OptionArmChance ---
$ws.Range("E33").Value = "25f"

# --- Row 34: OptionArmDamage ---
$ws.Range("E34").Value = "0.5f"

# --- Row 35: OptionArmDuration ---
$ws.Range("E35").Value = "3.0f"

# --- Row 37: OptionLegChance ---
$ws.Range("E37").Value = "30f"

# --- Row 38: OptionLegDamage ---
$ws.Range("E38").Value = "0.75f"

# --- Row 39: OptionLegDuration ---
$ws.Range("E39").Value = "3.5f"

# --- Row 41: OptionDismembermentChance ---
$ws.Range("E41").Value = "80f"

# --- Row 42: OptionDismembermentDamage ---
$ws.Range("E42").Value = "3.0f"

# --- Row 43: OptionDismembermentDuration ---
$ws.Range("E43").Value = "8.0f"
